# "save setting for deafult whastapp and default message when sending bill"
#
# Sprint 6 status updates on Sheet1 (Sprint Stories table):
#   C19 "I want to have place to store configurations: default message
#        when sending bill, default app to send bill"  -> DONE
#   C20 "Fix request permission problem"                -> DONE
#   C21 "I want to include other options to send my bill rather than
#        Whatsapp"                                      -> IN PROGRESS

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells that already carry the exact cell format (fill + border)
# used for each status so the bordered "Good"/"Neutral" table styles are
# reproduced exactly (plain Style-name assignment drops the border).
$refDone = $ws.Range("C2")        # "DONE" cell style (green, bordered)
$refInProgress = $ws.Range("A19") # "IN PROGRESS"/Sprint-label style (bordered)

$c19 = $ws.Range("C19")
$c19.Value = "DONE"
$refDone.Copy()
$c19.PasteSpecial(-4122)  # xlPasteFormats

$c20 = $ws.Range("C20")
$c20.Value = "DONE"
$refDone.Copy()
$c20.PasteSpecial(-4122)  # xlPasteFormats

$c21 = $ws.Range("C21")
$c21.Value = "IN PROGRESS"
$refInProgress.Copy()
$c21.PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Match the author's final selection in the saved file.
$null = $ws.Range("E17").Select()

$wb.Save()
